# Apply the workbook edit: remove the "exchangeable sodium percentage"
# indicator group from the mean_effect and TWW sheets, renumber the
# remaining sequence index on mean_effect, widen column A on TWW, and
# leave the selection/active-sheet state matching the saved workbook
# (Soil tab active).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("mean_effect")
$ws2 = $wb.Worksheets.Item("TWW")
$ws3 = $wb.Worksheets.Item("Soil")

# --- mean_effect: delete the "exchangeable sodium percentage" row (row 4)
$ws1.Range("A4").EntireRow.Delete()

# Renumber the sequence column (B) for the rows that shifted up so the
# index stays contiguous (3..8 instead of 4..9).
$ws1.Cells.Item(4, 2).Value = 3
$ws1.Cells.Item(5, 2).Value = 4
$ws1.Cells.Item(6, 2).Value = 5
$ws1.Cells.Item(7, 2).Value = 6
$ws1.Cells.Item(8, 2).Value = 7
$ws1.Cells.Item(9, 2).Value = 8

# --- TWW: delete the 5-row "exchangeable sodium percentage" block (rows 12-16)
$ws2.Range("A12:A16").EntireRow.Delete()

# Widen column A to fit the longer indicator labels.
$ws2.Columns.Item(1).ColumnWidth = 33.1667

# --- Selection / active sheet bookkeeping, matching the saved view state.
$ws1.Activate() | Out-Null
$ws1.Range("C3:I3").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("D7").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("D7:J10").Select() | Out-Null
